# "M16 Froze Token Embeddings + Decoder 12"
# Updates the per-epoch validation-accuracy values in column B (rows 3-115)
# to the newly recorded run, and moves the view/selection down to the
# bottom of the table where the run finished (row 115).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (column B, "Epoch Accuracy") whose value changed in the new run.
$rows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 22, 24, 25, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 105, 106, 107, 108, 109, 111, 112, 113, 115)

# Corresponding new accuracy values, in the same order as $rows.
$newValues = @(0.6875, 0.59375, 0.546875, 0.5, 0.453125, 0.34375, 0.375, 0.390625, 0.359375, 0.359375, 0.359375, 0.359375, 0.34375, 0.375, 0.375, 0.3125, 0.359375, 0.40625, 0.375, 0.375, 0.359375, 0.3125, 0.390625, 0.359375, 0.375, 0.375, 0.34375, 0.3125, 0.3125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.28125, 0.296875, 0.296875, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.3125, 0.515625, 0.40625, 0.546875, 0.53125, 0.34375, 0.484375, 0.46875, 0.421875, 0.4545454545454545)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $newValues[$i]
}

# The workbook was left scrolled near the bottom of the data with the
# last accuracy cell active.
$ws.Range("B115").Select()
